$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Agosto de 2020 a las 05:38"

# --- Swap Montserrat / Islas Malvinas rows (213 / 214) ---
# Row 213 held "Montserrat" data, row 214 held "Islas Malvinas" data.
# After the edit, row 213 shows "Islas Malvinas" and row 214 shows "Montserrat",
# with each country's stats travelling with its (now relocated) name.
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# --- Update country statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Muertes hoy, Muertes) ---

# Row 10: Peru
$ws.Range("B10").Value = 489680
$ws.Range("D10").Value = 335756
$ws.Range("E10").Value = 132423
$ws.Range("H10").Value = 21501

# Row 29: Kazajistan
$ws.Range("B29").Value = 100855
$ws.Range("C29").Value = 691
$ws.Range("D29").Value = 74677
$ws.Range("E29").Value = 24909

# Row 51: Honduras
$ws.Range("B51").Value = 48403
$ws.Range("C51").Value = 531
$ws.Range("D51").Value = 6805
$ws.Range("E51").Value = 40083
$ws.Range("G51").Value = 9
$ws.Range("H51").Value = 1515

# Row 154: Jamaica
$ws.Range("B154").Value = 1047
$ws.Range("C154").Value = 16
$ws.Range("D154").Value = 753
$ws.Range("E154").Value = 280

# Row 177: Mongolia
$ws.Range("D177").Value = 269
$ws.Range("E177").Value = 24

# Row 179: Camboya
$ws.Range("B179").Value = 268
$ws.Range("C179").Value = 2
$ws.Range("E179").Value = 48

# Row 213: now "Islas Malvinas" (was Montserrat's slot)
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

# Row 214: now "Montserrat" (was Islas Malvinas' slot)
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
